# Insert a new weekly price-report row for Cilantro at row 167, shifting the
# existing rows (old 167..228) down to (168..229). Mirrors the commit
# "Fruta / hortaliza, semanal" which adds one more weekly observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing data (rows 167-228) down by one row.
$ws.Rows.Item(167).Insert()

# Populate the newly inserted row 167 with the new weekly observation.
$ws.Cells.Item(167, 1).Value  = 10
$ws.Cells.Item(167, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(167, 3).Value  = "La Araucanía"
$ws.Cells.Item(167, 4).Value  = 44468
$ws.Cells.Item(167, 5).Value  = 9
$ws.Cells.Item(167, 6).Value  = 100112040
$ws.Cells.Item(167, 7).Value  = "Cilantro"
$ws.Cells.Item(167, 8).Value  = "Sin especificar"
$ws.Cells.Item(167, 9).Value  = "Primera"
$ws.Cells.Item(167, 10).Value = 30
$ws.Cells.Item(167, 11).Value = 4000
$ws.Cells.Item(167, 12).Value = 4000
$ws.Cells.Item(167, 13).Value = 4000
$ws.Cells.Item(167, 14).Value = "`$/docena de atados (2 kilos)"
$ws.Cells.Item(167, 15).Value = "Región Metropolitana"
$ws.Cells.Item(167, 16).Value = 2000
$ws.Cells.Item(167, 17).Value = 2
$ws.Cells.Item(167, 18).Value = "Hortaliza"
